$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - match formatting of the existing header cells (bold,
# bordered, centered/top aligned). Copy/PasteSpecial formats from the
# neighboring header cell (E1) so the existing style is reused rather than
# a brand-new (duplicate) style record being created.
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)   # xlPasteFormats

# Data cells F2:F11 with timestamps matching the diff
$timestamps = @(
    "2021-10-05 13:40:52.587925",
    "2021-10-05 13:40:52.587937",
    "2021-10-05 13:40:52.587941",
    "2021-10-05 13:40:52.587944",
    "2021-10-05 13:40:52.587948",
    "2021-10-05 13:40:52.587951",
    "2021-10-05 13:40:52.587954",
    "2021-10-05 13:40:52.587957",
    "2021-10-05 13:40:52.587961",
    "2021-10-05 13:40:52.587964"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
